$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = 0
$ws.Range("H132").Value = 2244163.2
$ws.Range("J132").Value = 6000
$ws.Range("L132").Value = 18000
$ws.Range("N132").Value = -23060
$ws.Range("H138").Value = 2004.8368
$ws.Range("I138").Value = 1065.1714
$ws.Range("K138").Value = 3195.5142
$ws.Range("M138").Value = 1944.4858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3429.9807
$ws.Range("I61").Value = 1127.6857
$ws.Range("K61").Value = 1127.6857
$ws.Range("M61").Value = -915.6857
$ws.Range("H132").Value = 1149.421
$ws.Range("I132").Value = 1168.8334
$ws.Range("K132").Value = 3506.5002
$ws.Range("M132").Value = -976.5001999999999
$ws.Range("H136").Value = 3429.9807
$ws.Range("I136").Value = 1127.6857
$ws.Range("K136").Value = 3383.0571
$ws.Range("M136").Value = -833.0571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 10558.546
$ws.Range("I20").Value = 13823.512
$ws.Range("K20").Value = 13823.512
$ws.Range("M20").Value = -13576.512
$ws.Range("H99").Value = 2286.4
$ws.Range("J99").Value = 2286.4
$ws.Range("L99").Value = 2286.4
$ws.Range("N99").Value = -5282.4
$ws.Range("H105").Value = 2500
$ws.Range("I105").Value = 2500
$ws.Range("K105").Value = 2500
$ws.Range("M105").Value = -753
$ws.Range("H106").Value = 11685.333
$ws.Range("J106").Value = 11685.333
$ws.Range("L106").Value = 11685.333
$ws.Range("N106").Value = -14209.333
$ws.Range("H134").Value = 1527.836
$ws.Range("I134").Value = 1230.3726
$ws.Range("J134").Value = 3044.9
$ws.Range("K134").Value = 3691.1178
$ws.Range("L134").Value = 9134.700000000001
$ws.Range("M134").Value = -1156.1178
$ws.Range("N134").Value = -14204.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 37000
$ws.Range("J51").Value = 37000
$ws.Range("L51").Value = 37000
$ws.Range("N51").Value = -38472
$ws.Range("H58").Value = 7288.0566
$ws.Range("I58").Value = 1219.7
$ws.Range("J58").Value = 25959.924
$ws.Range("K58").Value = 1219.7
$ws.Range("L58").Value = 25959.924
$ws.Range("M58").Value = -1016.7
$ws.Range("N58").Value = -26365.924
$ws.Range("H61").Value = 37000
$ws.Range("J61").Value = 37000
$ws.Range("L61").Value = 37000
$ws.Range("N61").Value = -37696
$ws.Range("H99").Value = 17798.2
$ws.Range("I99").Value = 15667.333
$ws.Range("K99").Value = 15667.333
$ws.Range("M99").Value = -14169.333
$ws.Range("H122").Value = 3449.5
$ws.Range("J122").Value = 3450
$ws.Range("L122").Value = 10350
$ws.Range("N122").Value = -15250
$ws.Range("H126").Value = 17798.2
$ws.Range("I126").Value = 15667.333
$ws.Range("K126").Value = 47001.999
$ws.Range("M126").Value = -44531.999
$ws.Range("H136").Value = 7288.0566
$ws.Range("I136").Value = 1219.7
$ws.Range("J136").Value = 25959.924
$ws.Range("K136").Value = 3659.1
$ws.Range("L136").Value = 77879.772
$ws.Range("M136").Value = -1109.1
$ws.Range("N136").Value = -82979.772

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2948.75
$ws.Range("I5").Value = 2948.75
$ws.Range("K5").Value = 8846.25
$ws.Range("M5").Value = -8734.25
$ws.Range("H69").Value = 4803.4346
$ws.Range("I69").Value = 2740
$ws.Range("K69").Value = 8220
$ws.Range("M69").Value = -7409
$ws.Range("H72").Value = 4803.4346
$ws.Range("I72").Value = 2740
$ws.Range("K72").Value = 24660
$ws.Range("M72").Value = -20604
$ws.Range("H126").Value = 5000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 0
$ws.Range("L126").ClearContents()
$ws.Range("M126").Value = 15000
$ws.Range("N126").Value = -24880
$ws.Range("H129").Value = 2982.125
$ws.Range("I129").Value = 2029.5
$ws.Range("K129").Value = 6088.5
$ws.Range("M129").Value = -1088.5
$ws.Range("H132").Value = 1699
$ws.Range("I132").Value = 998
$ws.Range("K132").Value = 8982
$ws.Range("M132").Value = -6452
$ws.Range("H133").Value = 6419.5835
$ws.Range("I133").Value = 3397.25
$ws.Range("J133").Value = 7930.75
$ws.Range("K133").Value = 10191.75
$ws.Range("L133").Value = 23792.25
$ws.Range("M133").Value = -5131.75
$ws.Range("N133").Value = -33912.25
$ws.Range("H134").Value = 719.5
$ws.Range("I134").Value = 719.5
$ws.Range("K134").Value = 2158.5
$ws.Range("M134").Value = 2911.5
$ws.Range("H135").Value = 2948.75
$ws.Range("I135").Value = 2948.75
$ws.Range("K135").Value = 26538.75
$ws.Range("M135").Value = -24003.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3237.5625
$ws.Range("I132").Value = 3248.6296
$ws.Range("J132").Value = 3177.8
$ws.Range("K132").Value = 9745.888800000001
$ws.Range("L132").Value = 9533.400000000001
$ws.Range("M132").Value = -7215.888800000001
$ws.Range("N132").Value = -14593.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3315.4
$ws.Range("I7").Value = 2730.8
$ws.Range("J7").Value = 3900
$ws.Range("K7").Value = 2730.8
$ws.Range("L7").Value = 3900
$ws.Range("M7").Value = -2618.8
$ws.Range("N7").Value = -4124
$ws.Range("H16").Value = 1925.7858
$ws.Range("I16").Value = 1632.7273
$ws.Range("K16").Value = 1632.7273
$ws.Range("M16").Value = -1462.7273
$ws.Range("H40").Value = 3754.2222
$ws.Range("I40").Value = 3113
$ws.Range("K40").Value = 3113
$ws.Range("M40").Value = -2977
$ws.Range("H126").Value = 3315.4
$ws.Range("I126").Value = 2730.8
$ws.Range("J126").Value = 3900
$ws.Range("K126").Value = 8192.400000000001
$ws.Range("L126").Value = 11700
$ws.Range("M126").Value = -5722.400000000001
$ws.Range("N126").Value = -16640
$ws.Range("H134").Value = 47722.285
$ws.Range("J134").Value = 47722.285
$ws.Range("L134").Value = 47722.285
$ws.Range("N134").Value = -57862.285

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1461.1111
$ws.Range("I113").Value = 967.7143
$ws.Range("K113").Value = 2903.1429
$ws.Range("M113").Value = -733.1428999999998
$ws.Range("H126").Value = 280231
$ws.Range("I126").Value = 1609.9333
$ws.Range("K126").Value = 4829.7999
$ws.Range("M126").Value = -2359.7999
$ws.Range("H136").Value = 9279.493
$ws.Range("I136").Value = 10336.627
$ws.Range("K136").Value = 31009.881
$ws.Range("M136").Value = -28459.881
